# Update NATMI TPM-derived receptor expression values (Gdf5-Bmpr1a, row 2)
# and the resulting recalculated specificity / edge-weight columns for rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - new receptor average / total expression values (source TPM update)
$ws.Range("M2").Value = 2.341355666666667
$ws.Range("N2").Value = 7.024067000000001

# Row 2 - recalculated derived-specificity / edge-weight columns
$ws.Range("O2").Value = 0.03973512964576821
$ws.Range("P2").Value = 0.0397351296457682
$ws.Range("Q2").Value = 0.07118657768933334
$ws.Range("R2").Value = 0.6406791992040001
$ws.Range("S2").Value = 0.03973512964576821
$ws.Range("T2").Value = 0.0397351296457682

# Row 3 - recalculated derived-specificity columns (M3/N3 themselves unchanged)
$ws.Range("O3").Value = 0.5779093692199981
$ws.Range("P3").Value = 0.5779093692199981
$ws.Range("Q3").Value = 1.035340530561333
$ws.Range("S3").Value = 0.5779093692199981
$ws.Range("T3").Value = 0.5779093692199981

# Row 4 - recalculated derived-specificity columns (M4/N4 themselves unchanged)
$ws.Range("O4").Value = 0.3823555011342337
$ws.Range("P4").Value = 0.3823555011342337
$ws.Range("S4").Value = 0.3823555011342337
$ws.Range("T4").Value = 0.3823555011342337
